$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "1.3.1.1.f" -> "1.3.1.1f" typo in the Russian header (column B, row 1) ---
$ws.Range("B1").Value = "1.3.1.1f Доля лиц, получающих пенсии и пособия по инвалидности к общей численности населения"

# --- Add a new year column (T) for 2023, copying the formatting from column S ---
$ws.Range("S2:S5").Copy($ws.Range("T2:T5"))

# T2 has no value (just formatting/border), leave it blank.
$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 217222
$ws.Range("T5").Value = 2.9794303052841493

# --- Row heights for rows 2 and 4 changed slightly ---
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 29.25
